$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header strings: ctb_banddis -> comparison_df, frs_banddis -> V2
$ws.Range("B1").Value = "comparison_df"
$ws.Range("C1").Value = "V2"

# Updated C (comparison_df) and D (banddis_diff) values, rows 2-9
$ws.Range("C2").Value = 20.646515533165406
$ws.Range("D2").Value = 3.294165716876961

$ws.Range("C3").Value = 18.564231738035264
$ws.Range("D3").Value = 0.9578499221717287

$ws.Range("C4").Value = 22.55247691015953
$ws.Range("D4").Value = -0.6596654311633507

$ws.Range("C5").Value = 17.094878253568428
$ws.Range("D5").Value = -1.4874216504614353

$ws.Range("C6").Value = 10.831234256926953
$ws.Range("D6").Value = -1.0900364577650663

$ws.Range("C7").Value = 6.095717884130982
$ws.Range("D7").Value = -0.9244230334995907

$ws.Range("C8").Value = 3.7615449202350963
$ws.Range("D8").Value = -0.23176881804699567

$ws.Range("C9").Value = 0.4534005037783375
$ws.Range("D9").Value = 0.14129975188774868
